# "added music + difficulty"
# Fill in missing LICENSE TYPE values, add two JUNGLE rows' data to the
# previously-blank rows 6 and 7, and append two new HOUSE rows (8 and 9)
# to the TRACK ID LIST table.

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Sets a cell's paragraph content from a fragment of run-level WordprocessingML
# (runs / proofErr markers, etc.) and makes sure the cell ends up with exactly
# one paragraph (InsertXML can leave a stray empty leading paragraph behind on
# freshly-added rows, so clean that up afterwards).
function Set-CellRuns($cell, [string]$innerXml) {
    $xml = "<w:p $wNs>$innerXml</w:p>"
    $cell.Range.InsertXML($xml)
    while ($cell.Range.Paragraphs.Count -gt 1) {
        $cell.Range.Paragraphs.Item(1).Range.Delete()
    }
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 2 (ID 1, Techno001): LICENSE TYPE was blank -> "Royalty free"
$t.Cell(2, 5).Range.Text = "Royalty free"

# --- Row 3 (ID 2, Techno002): LICENSE TYPE "Owned by me" -> "Original content"
$t.Cell(3, 5).Range.Text = "Original content"

# --- Row 4 (ID 3, Techno003): LICENSE TYPE was blank -> "Royalty free"
$t.Cell(4, 5).Range.Text = "Royalty free"

# --- Row 5 (ID 4, Techno004): LICENSE TYPE was blank -> "Royalty free"
$t.Cell(5, 5).Range.Text = "Royalty free"

# --- Row 6 (ID 5, House001): LICENSE TYPE "Owned by me" -> "Original content"
$t.Cell(6, 5).Range.Text = "Original content"

# --- Row 7 (previously fully blank) -> ID 6 / Jungle001 / JUNGLE / Sappheiros - Universe / Creative Commons (BY 3.0)
$t.Cell(7, 1).Range.Text = "6"
$t.Cell(7, 2).Range.Text = "Jungle001"
$t.Cell(7, 3).Range.Text = "JUNGLE"
$inner7 = '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>Sappheiros</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> - Universe</w:t></w:r>'
Set-CellRuns $t.Cell(7, 4) $inner7
$t.Cell(7, 5).Range.Text = "Creative Commons (BY 3.0)"

# --- Row 8 (previously fully blank) -> ID 7 / Jungle002 / JUNGLE / Ripple (en dash) On Your Mind / NCS (en dash) Copyright free
$t.Cell(8, 1).Range.Text = "7"
$t.Cell(8, 2).Range.Text = "Jungle002"
$t.Cell(8, 3).Range.Text = "JUNGLE"
$t.Cell(8, 4).Range.Text = [string]::Format("Ripple {0} On Your Mind", [char]0x2013)
$t.Cell(8, 5).Range.Text = [string]::Format("NCS {0} Copyright free", [char]0x2013)

# --- New row 9 -> ID 8 / House002 / HOUSE / Vidaloca - Housess / Royalty free
$row9 = $t.Rows.Add()
$row9.Cells(1).Range.Text = "8"
$row9.Cells(2).Range.Text = "House002"
$row9.Cells(3).Range.Text = "HOUSE"
$inner9 = '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>Vidaloca</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> - </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>Housess</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>'
Set-CellRuns $row9.Cells(4) $inner9
$row9.Cells(5).Range.Text = "Royalty free"

# --- New row 10 -> ID 9 / House003 (two runs) / HOUSE /
#     Jorja Smith x Preditah - On My Mind (Imacci Remix) / Original content
$row10 = $t.Rows.Add()
$row10.Cells(1).Range.Text = "9"
$inner10b = '<w:r><w:t>House00</w:t></w:r>' `
    + '<w:r><w:t>3</w:t></w:r>'
Set-CellRuns $row10.Cells(2) $inner10b
$row10.Cells(3).Range.Text = "HOUSE"
$inner10d = '<w:r><w:t xml:space="preserve">Jorja Smith x </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>Preditah</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> - On My Mind (</w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>Imacci</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> Remix)</w:t></w:r>'
Set-CellRuns $row10.Cells(4) $inner10d
$row10.Cells(5).Range.Text = "Original content"
